# Renames a handful of pseudo-code identifiers inside the "for loop"
# code sample textbox (TextBox 15, inside Group 11) that is repeated,
# essentially identically, on every slide of the deck:
#   Write(            -> write(
#   ToInt32(          -> to_integer(
#   ReadLine          -> read_line
#   WriteLine(        -> write_line(
#
# The COM text model here round-trips the run text using .Characters()
# sub-ranges; to avoid PowerPoint silently splitting a single <a:r> run
# into several runs, each replacement below always spans the *entire*
# original run's text (never just the changed substring).
#
# Note: reading the "…" (U+2026) ellipsis character back out through
# TextRange.Text/.Characters(...).Text is lossy in this runtime (it
# comes back as "."), so searches use an ASCII-safe prefix/suffix and
# the replacement text is always supplied as a fresh literal (never
# built from a read-back value).

function Replace-RunSpan {
    param(
        $TextRange,
        [string]$FindPrefix,
        [int]$OldLength,
        [string]$NewText
    )

    $full = $TextRange.Text
    $idx = $full.IndexOf($FindPrefix)
    if ($idx -lt 0) {
        Write-Host "WARNING: prefix not found: [$FindPrefix]"
        return
    }

    $startChar = $idx + 1   # PowerPoint Characters() is 1-based
    $span = $TextRange.Characters($startChar, $OldLength)
    $span.Text = $NewText
}

function Update-ForLoopCodeBox {
    param($TextRange)

    # Para 1: Write("Count to: ");
    Replace-RunSpan $TextRange 'Write("Count to: ");' 20 'write("Count to: ");'

    # Para 2, run 1: target = ToInt32(
    Replace-RunSpan $TextRange 'target = ToInt32(' 17 'target = to_integer('

    # Para 2, run 2: ReadLine
    Replace-RunSpan $TextRange 'ReadLine' 8 'read_line'

    # Para 5, run 1: "    WriteLine(" (four leading spaces)
    Replace-RunSpan $TextRange '    WriteLine(' 14 '    write_line('

    # Para 7: WriteLine("Bye…");
    Replace-RunSpan $TextRange 'WriteLine("Bye' 18 'write_line("Bye…");'
}

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.Name -eq "Group 11" -and $shape.Type -eq 6) {
            for ($k = 1; $k -le $shape.GroupItems.Count; $k++) {
                $item = $shape.GroupItems.Item($k)
                if ($item.Name -eq "TextBox 15") {
                    Update-ForLoopCodeBox $item.TextFrame.TextRange
                }
            }
        }
    }
}
